# Auto-update: dice_jobs_list.xlsx
# Appends three new job listing rows (83-85) to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{
        Title = "Hiring for a Golang Developer in McLean, VA / Richmond VA"
        URL = "https://www.dice.com/job-detail/fe5d1d61-23a5-416d-bbae-7368d6cb757f"
        Location = "Richmond, Virginia"
        Employment_Type = "Full-time, Contract"
        Salary = "Depends on Experience"
        Company = "InfiCare Technologies"
    },
    @{
        Title = "Golang Architect / Principal Backend Architect Only Local to GA"
        URL = "https://www.dice.com/job-detail/0e936e90-924d-4b67-b0c7-03e143cfbeb8"
        Location = "Atlanta, Georgia"
        Employment_Type = "Third Party"
        Salary = "Depends on Experience"
        Company = "Dahl Consulting"
    },
    @{
        Title = "Google CCAI Tech Lead"
        URL = "https://www.dice.com/job-detail/bf26aea7-da55-4c7b-8a1a-53e6f37a4dc0"
        Location = "Hybrid in Hartford, Connecticut"
        Employment_Type = "Contract, Third Party"
        Salary = "Depends on Experience"
        Company = "VDart, Inc."
    }
)

$startRow = 83
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row.Title
    $ws.Cells.Item($r, 2).Value = $row.URL
    $ws.Cells.Item($r, 3).Value = $row.Location
    $ws.Cells.Item($r, 4).Value = $row.Employment_Type
    $ws.Cells.Item($r, 5).Value = $row.Salary
    $ws.Cells.Item($r, 6).Value = $row.Company
}
